$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellDate($row, $col, $text) {
    $cell = $tbl.Rows.Item($row).Cells.Item($col)
    $rng = $cell.Range
    $rng.Text = $text
    $rng.Font.Name = "Times New Roman"
    $rng.Font.NameAscii = "Times New Roman"
    $rng.Font.Size = 14
}

# Row 6 (topic "Планування конструювання програмного забезпечення" / Л03):
#   ПЗ-42 "Заняття" cell (column 3) gets 14.09
Set-CellDate 6 3 "14.09"

# Row 7 (topic "Розробка технічного завдання на програмне забезпечення" / ПР01):
#   ПЗ-41 "Заняття" (col 1) = 13.09, ПЗ-41 "Здача" (col 2) = 17.09
#   ПЗ-42 "Заняття" (col 3) = 15.09, ПЗ-42 "Здача" (col 4) = 19.09
Set-CellDate 7 1 "13.09"
Set-CellDate 7 2 "17.09"
Set-CellDate 7 3 "15.09"
Set-CellDate 7 4 "19.09"

# Row 8 (topic "Планування конструювання програмного з забезпечення" / ЛР02):
#   ПЗ-41 "Заняття" (col 1) = 14.09, ПЗ-41 "Здача" (col 2) = 18.09
Set-CellDate 8 1 "14.09"
Set-CellDate 8 2 "18.09"

Write-Output "done"
